$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old SecurityComponent rows (rows 6 and 7) ---
# Their data (sec-identity / sec-threat) gets folded into rows 4 and 5 as
# BudgetCode/EBSCostCenter tag pairs, same shape as the cost- groups, so
# the table drops from 7 data rows down to 4.
$ws.Rows("6:7").Delete()

# --- Row 2: cost-aexeo (BudgetCode + EBSCostCenter) ---
$ws.Range("A2").Value = "TagGroup"
$ws.Range("C2").Value = "cost-aexeo"
$ws.Range("D2").Value = "BudgetCode"
$ws.Range("E2").Value = "PROJECT-1212121"
$ws.Range("F2").Value = "EBSCostCenter"
$ws.Range("G2").Value = "'4000009"

# --- Row 3: cost-wan (BudgetCode + EBSCostCenter) ---
$ws.Range("A3").Value = "TagGroup"
$ws.Range("C3").Value = "cost-wan"
$ws.Range("D3").Value = "BudgetCode"
$ws.Range("E3").Value = "PROJECT-232222"
$ws.Range("F3").Value = "EBSCostCenter"
$ws.Range("G3").Value = "'4000002"

# --- Row 4/5: sec-identity / sec-threat, now BudgetCode + EBSCostCenter pairs ---
$ws.Range("A4").Value = "TagGroup"
$ws.Range("C4").Value = "sec-identity"
$ws.Range("D4").Value = "BudgetCode"

$ws.Range("A5").Value = "TagGroup"
$ws.Range("C5").Value = "sec-threat"
$ws.Range("D5").Value = "BudgetCode"

# --- Header row: TagKey/TagValue -> TagKey1..3/TagValue1..3 (first new columns) ---
$ws.Range("D1").Value = "TagKey1"
$ws.Range("E1").Value = "TagValue1"
$ws.Range("F1").Value = "TagKey2"
$ws.Range("G1").Value = "TagValue2"
$ws.Range("H1").Value = "TagKey3"
$ws.Range("I1").Value = "TagValue3"

# --- Fill the new BudgetCode/EBSCostCenter values for the security rows ---
$ws.Range("E5").Value = "PROJECT-901233"
$ws.Range("E4").Value = "PROJECT-896775"
$ws.Range("G4").Value = "'4000033"
$ws.Range("G5").Value = "'4000055"
$ws.Range("F4").Value = "EBSCostCenter"
$ws.Range("F5").Value = "EBSCostCenter"

# --- Header row: add TagKey4/TagValue4/TagKey5/TagValue5 (room for API-limit paging) ---
$ws.Range("J1").Value = "TagKey4"
$ws.Range("K1").Value = "TagValue4"
$ws.Range("L1").Value = "TagKey5"
$ws.Range("M1").Value = "TagValue5"

# --- Re-establish the shared concatenation formula over the (now smaller) B2:B5 ---
$ws.Range("B2:B5").Formula = '=A2&"/"&C2&"/"&D2&"/"&E2'

# --- New column widths for the additional tag-pair columns, matching D:E ---
$ws.Columns("F").ColumnWidth = 18.166666666666668
$ws.Columns("G").ColumnWidth = 16.5
$ws.Columns("H").ColumnWidth = 18.166666666666668
$ws.Columns("I").ColumnWidth = 16.5

# --- View state: new cell selected after the edits ---
$ws.Range("M2").Select()

Write-Output "done"
